$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recompute NATMI ligand/receptor TPM-derived statistics and edge weights
# (columns G-T) for rows 2-10 using the updated TPM expression matrix.

# Row 2
$ws.Cells.Item(2, 7).Value = 0.2379653333333333
$ws.Cells.Item(2, 8).Value = 0.713896
$ws.Cells.Item(2, 9).Value = 0.0004000853538884766
$ws.Cells.Item(2, 10).Value = 0.0004000853538884766
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.292163333333333
$ws.Cells.Item(2, 14).Value = 3.87649
$ws.Cells.Item(2, 15).Value = 0.00934486532303975
$ws.Cells.Item(2, 16).Value = 0.00934486532303975
$ws.Cells.Item(2, 17).Value = 0.3074900783377778
$ws.Cells.Item(2, 18).Value = 2.76741070504
$ws.Cells.Item(2, 19).Value = 0.000003738743749808512
$ws.Cells.Item(2, 20).Value = 0.000003738743749808512

# Row 3
$ws.Cells.Item(3, 7).Value = 0.2379653333333333
$ws.Cells.Item(3, 8).Value = 0.713896
$ws.Cells.Item(3, 9).Value = 0.0004000853538884766
$ws.Cells.Item(3, 10).Value = 0.0004000853538884766
$ws.Cells.Item(3, 15).Value = 0.8196320797583818
$ws.Cells.Item(3, 16).Value = 0.8196320797583818
$ws.Cells.Item(3, 17).Value = 26.96975544331111
$ws.Cells.Item(3, 18).Value = 242.7277989898
$ws.Cells.Item(3, 19).Value = 0.0003279227906884803
$ws.Cells.Item(3, 20).Value = 0.0003279227906884803

# Row 4
$ws.Cells.Item(4, 7).Value = 0.2379653333333333
$ws.Cells.Item(4, 8).Value = 0.713896
$ws.Cells.Item(4, 9).Value = 0.0004000853538884766
$ws.Cells.Item(4, 10).Value = 0.0004000853538884766
$ws.Cells.Item(4, 13).Value = 23.648251
$ws.Cells.Item(4, 14).Value = 70.94475299999999
$ws.Cells.Item(4, 15).Value = 0.1710230549185785
$ws.Cells.Item(4, 16).Value = 0.1710230549185785
$ws.Cells.Item(4, 17).Value = 5.627463931965333
$ws.Cells.Item(4, 18).Value = 50.64717538768799
$ws.Cells.Item(4, 19).Value = 0.00006842381945018785
$ws.Cells.Item(4, 20).Value = 0.00006842381945018785

# Row 5
$ws.Cells.Item(5, 9).Value = 0.9842542228653065
$ws.Cells.Item(5, 10).Value = 0.9842542228653065
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.292163333333333
$ws.Cells.Item(5, 14).Value = 3.87649
$ws.Cells.Item(5, 15).Value = 0.00934486532303975
$ws.Cells.Item(5, 16).Value = 0.00934486532303975
$ws.Cells.Item(5, 17).Value = 756.4596033113089
$ws.Cells.Item(5, 18).Value = 6808.136429801781
$ws.Cells.Item(5, 19).Value = 0.009197723156309441
$ws.Cells.Item(5, 20).Value = 0.009197723156309441

# Row 6
$ws.Cells.Item(6, 9).Value = 0.9842542228653065
$ws.Cells.Item(6, 10).Value = 0.9842542228653065
$ws.Cells.Item(6, 15).Value = 0.8196320797583818
$ws.Cells.Item(6, 16).Value = 0.8196320797583818
$ws.Cells.Item(6, 19).Value = 0.806726335698061
$ws.Cells.Item(6, 20).Value = 0.806726335698061

# Row 7
$ws.Cells.Item(7, 9).Value = 0.9842542228653065
$ws.Cells.Item(7, 10).Value = 0.9842542228653065
$ws.Cells.Item(7, 13).Value = 23.648251
$ws.Cells.Item(7, 14).Value = 70.94475299999999
$ws.Cells.Item(7, 15).Value = 0.1710230549185785
$ws.Cells.Item(7, 16).Value = 0.1710230549185785
$ws.Cells.Item(7, 17).Value = 13844.18371036654
$ws.Cells.Item(7, 18).Value = 124597.6533932989
$ws.Cells.Item(7, 19).Value = 0.1683301640109361
$ws.Cells.Item(7, 20).Value = 0.1683301640109361

# Row 8
$ws.Cells.Item(8, 7).Value = 9.127409
$ws.Cells.Item(8, 8).Value = 27.382227
$ws.Cells.Item(8, 9).Value = 0.01534569178080505
$ws.Cells.Item(8, 10).Value = 0.01534569178080505
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.292163333333333
$ws.Cells.Item(8, 14).Value = 3.87649
$ws.Cells.Item(8, 15).Value = 0.00934486532303975
$ws.Cells.Item(8, 16).Value = 0.00934486532303975
$ws.Cells.Item(8, 17).Value = 11.79410323813667
$ws.Cells.Item(8, 18).Value = 106.14692914323
$ws.Cells.Item(8, 19).Value = 0.0001434034229805012
$ws.Cells.Item(8, 20).Value = 0.0001434034229805012

# Row 9
$ws.Cells.Item(9, 7).Value = 9.127409
$ws.Cells.Item(9, 8).Value = 27.382227
$ws.Cells.Item(9, 9).Value = 0.01534569178080505
$ws.Cells.Item(9, 10).Value = 0.01534569178080505
$ws.Cells.Item(9, 15).Value = 0.8196320797583818
$ws.Cells.Item(9, 16).Value = 0.8196320797583818
$ws.Cells.Item(9, 17).Value = 1034.453149594942
$ws.Cells.Item(9, 18).Value = 9310.078346354474
$ws.Cells.Item(9, 19).Value = 0.01257782126963235
$ws.Cells.Item(9, 20).Value = 0.01257782126963235

# Row 10
$ws.Cells.Item(10, 7).Value = 9.127409
$ws.Cells.Item(10, 8).Value = 27.382227
$ws.Cells.Item(10, 9).Value = 0.01534569178080505
$ws.Cells.Item(10, 10).Value = 0.01534569178080505
$ws.Cells.Item(10, 13).Value = 23.648251
$ws.Cells.Item(10, 14).Value = 70.94475299999999
$ws.Cells.Item(10, 15).Value = 0.1710230549185785
$ws.Cells.Item(10, 16).Value = 0.1710230549185785
$ws.Cells.Item(10, 17).Value = 215.847259011659
$ws.Cells.Item(10, 18).Value = 1942.625331104931
$ws.Cells.Item(10, 19).Value = 0.0026244670881922
$ws.Cells.Item(10, 20).Value = 0.0026244670881922
